$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest crypto data.
# Leading apostrophe on D-column writes keeps the cell as Text (matches
# the source data's formatting, e.g. trailing zeros like '8.30' and
# multi-dot values like '62.078.71' which must not be coerced to numbers).

$ws.Range("D2").Value = "'62.078.71"
$ws.Range("E2").Value = "  +2.86%  "

$ws.Range("D3").Value = "'2.418.72"
$ws.Range("E3").Value = "  +4.12%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'559.46"
$ws.Range("E5").Value = "  +2.64%  "

$ws.Range("D6").Value = "'138.89"
$ws.Range("E6").Value = "  +6.04%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "'0.583"
$ws.Range("E8").Value = "  +0.81%  "

$ws.Range("D9").Value = "'2.416.13"
$ws.Range("E9").Value = "  +4.19%  "

$ws.Range("E10").Value = "  +3.21%  "

$ws.Range("E11").Value = "  +3.95%  "

$ws.Range("E12").Value = "  -0.03%  "

$ws.Range("D14").Value = "'25.86"
$ws.Range("E14").Value = "  +9.58%  "

$ws.Range("D15").Value = "'2.847.65"
$ws.Range("E15").Value = "  +4.12%  "

$ws.Range("D16").Value = "'62.010.49"
$ws.Range("E16").Value = "  +2.79%  "

$ws.Range("E17").Value = "  +5.02%  "

$ws.Range("D18").Value = "'2.421.89"
$ws.Range("E18").Value = "  +3.98%  "

$ws.Range("D19").Value = "'11.09"
$ws.Range("E19").Value = "  +4.92%  "

$ws.Range("D20").Value = "'343.63"
$ws.Range("E20").Value = "  +9.44%  "

$ws.Range("E21").Value = "  +2.32%  "

$ws.Range("D22").Value = "'6.85"
$ws.Range("E22").Value = "  +3.03%  "

$ws.Range("E23").Value = "  +0.27%  "

$ws.Range("D24").Value = "'64.97"

$ws.Range("D25").Value = "'0.172"
$ws.Range("E25").Value = "  -0.71%  "

$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.57%  "

$ws.Range("D27").Value = "'8.30"
$ws.Range("E27").Value = "  +5.67%  "

$ws.Range("E28").Value = "  +10.88%  "

$ws.Range("E29").Value = "  +14.47%  "

$ws.Range("D30").Value = "'0.0₃0781"
$ws.Range("E30").Value = "  +6.80%  "

$ws.Range("E31").Value = "  +3.81%  "

$ws.Range("E32").Value = "  -1.02%  "

$ws.Range("D33").Value = "'6.32"
$ws.Range("E33").Value = "  +6.40%  "

$ws.Range("E34").Value = "  +3.14%  "

$ws.Range("D35").Value = "'0.395"
$ws.Range("E35").Value = "  +3.78%  "

$ws.Range("D36").Value = "'376.76"
$ws.Range("E36").Value = "  +16.67%  "

$ws.Range("D37").Value = "'18.52"
$ws.Range("E37").Value = "  +4.02%  "

$ws.Range("D38").Value = "'4.48"
$ws.Range("E38").Value = "  +10.42%  "

$ws.Range("E39").Value = "  -0.01%  "

$ws.Range("E40").Value = "  +0.09%  "

$ws.Range("E41").Value = "  +8.51%  "

$ws.Range("D42").Value = "'39.08"
$ws.Range("E42").Value = "  +3.04%  "

$ws.Range("D43").Value = "'145.46"
$ws.Range("E43").Value = "  +5.44%  "

$ws.Range("D44").Value = "'3.66"
$ws.Range("E44").Value = "  +4.90%  "

$ws.Range("D45").Value = "'20.63"
$ws.Range("E45").Value = "  +8.27%  "

$ws.Range("E46").Value = "  +1.87%  "

$ws.Range("D48").Value = "'0.0518"
$ws.Range("E48").Value = "  +4.81%  "

$ws.Range("D49").Value = "'18.02"
$ws.Range("E49").Value = "  +6.61%  "

$ws.Range("E50").Value = "  +3.22%  "

$ws.Range("D51").Value = "'0.0₆0223"
$ws.Range("E51").Value = "  +4.36%  "
